# Fruta / hortaliza, semanal
# Insert a new weekly data row (row 45) for "Vega Monumental Concepción - Ciruela",
# shifting the existing rows 45-62 down to 46-63, and update row 44 with the
# latest week's corrected price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 45 (pushes old rows 45..62 down to 46..63).
$ws.Rows.Item(45).Insert()

# 2) Update row 44 (the newest entry) with its corrected values.
$ws.Range("D44").Value = 44642
$ws.Range("O44").Value = 9000
$ws.Range("P44").Value = 8455
$ws.Range("S44").Value = 470

# 3) Populate the newly inserted row 45 with the data that used to be in row 44
#    (same market record, prior week), now re-filed one row down.
$ws.Range("A45").Value = 11
$ws.Range("B45").Value = "Vega Monumental Concepción"
$ws.Range("C45").Value = "Bíobío"
$ws.Range("D45").Value = 44637
$ws.Range("E45").Value = 8
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100103
$ws.Range("H45").Value = "Frutos de hueso (carozo)"
$ws.Range("I45").Value = 100103002
$ws.Range("J45").Value = "Ciruela"
$ws.Range("K45").Value = "Angeleno"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 220
$ws.Range("N45").Value = 8000
$ws.Range("O45").Value = 8500
$ws.Range("P45").Value = 8273
$ws.Range("Q45").Value = "$/bandeja 18 kilos granel"
$ws.Range("R45").Value = "Provincia de Curicó"
$ws.Range("S45").Value = 460
$ws.Range("T45").Value = 18
